# Timesheet update: fill in Mohammed's (col F) hours for the week of
# 2017-03-13 .. 2017-03-17, and Joel's (col E) hours for 2017-03-17/18,
# then leave the selection where the author left off editing (E10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = 9
$ws.Range("F5").Value = 7
$ws.Range("F6").Value = 10
$ws.Range("F7").Value = 12

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 6

$ws.Range("E9").Value = 7
$ws.Range("F9").Value = 8

$ws.Range("E10").Select()
